$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$fallbackBody = "Beste collega,`nOnderstaande e-mail kon niet automatisch worden beantwoord door het AI-systeem. Wil je deze even opvolgen?`n📩 Originele afzender: klantenservice@testbedrijf123.nl`n📝 Onderwerp: Demo inplannen`n🔎 Reden: Interne afzender`n━━━━━━━━━━━━━━━━━━━━━━━━━━━`n✉️ Bericht:`nKun je vrijdag om 11:00 een demo inplannen bij Van Dijk?`n━━━━━━━━━━━━━━━━━━━━━━━━━━━`nMet vriendelijke groet,`nMailMind Automatische Assistent`n—`n[Bedrijfsnaam]`nklantenservice@bedrijf.nl`nwww.bedrijf.nl"

$rows = @(
    ,@("Re: Demo inplannen", "admin@testbedrijf123.nl", "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl.", "Intern verzoek / Actie voor medewerker", "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl.", "2025-08-13 21:46:07", "Nee", "Ja", "Nee", "Nee")
    ,@("[Fallback] Handmatige opvolging: Demo inplannen", "admin@testbedrijf123.nl", $fallbackBody, "Intern verzoek / Actie voor medewerker", "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl.", "2025-08-13 21:46:08", "Nee", "Ja", "Nee", "Nee")
    ,@("Demo inplannen", "klantenservice@testbedrijf123.nl", "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?", "Intern verzoek / Actie voor medewerker", "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl.", "2025-08-13 21:46:09", "Nee", "Ja", "Nee", "Nee")
    ,@("Demo inplannen", "klantenservice@testbedrijf123.nl", "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?", "Intern verzoek / Actie voor medewerker", "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl.", "2025-08-13 21:46:35", "Nee", "Ja", "Nee", "Nee")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Extend conditional formatting ranges from row 13 to row 17 for D, G, H, I, J
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $ws.Range("$($col)2:$($col)13")
    $newRange = $ws.Range("$($col)2:$($col)17")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update Dashboard count
$dash.Range("B2").Value = 16
